# basic_examples.xlsx edit
#
# Adds three new example worksheets ("table with dividers",
# "data with merged cells", "merged cells header"), adds the shared string
# "score", and updates the view state (selection / active tab) of a couple
# of the pre-existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Create the three new worksheets (appended after the existing ones, in
# order) first, so tab order / sheetId / r:id all come out right.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "table with dividers"

$ws5 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws5.Name = "data with merged cells"

$ws6 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws6.Name = "merged cells header"

# ---------------------------------------------------------------------------
# Sheet 6 first: establishes the "bold + centered" cell style used by its
# merged "score" header cell.
# ---------------------------------------------------------------------------
$ws6.Columns.Item(1).ColumnWidth = 16.25

$ws6.Range("A1").Value = "home team"
$ws6.Range("B1").Value = "guest team"
$ws6.Range("C1").Value = "score"
$ws6.Range("E1").Value = "date"
$ws6.Range("A1:E1").Font.Bold = $true
$ws6.Range("C1:D1").Merge()
$ws6.Range("C1:D1").HorizontalAlignment = -4108

$ws6.Range("A2").Value = "Manchester City"
$ws6.Range("B2").Value = "RB Leipzig"
$ws6.Range("C2").Value = 6
$ws6.Range("D2").Value = 3
$ws6.Range("E2").Value = 44454
$ws6.Range("E2").NumberFormat = "d-mmm"

$ws6.Range("A3").Value = "Club Brugge"
$ws6.Range("B3").Value = "PSG"
$ws6.Range("C3").Value = 1
$ws6.Range("D3").Value = 1
$ws6.Range("E3").Value = 44454
$ws6.Range("E3").NumberFormat = "d-mmm"

$ws6.Range("A4").Value = "RB Leipzig"
$ws6.Range("B4").Value = "Club Brugge"
$ws6.Range("C4").Value = 1
$ws6.Range("D4").Value = 2
$ws6.Range("E4").Value = 44467
$ws6.Range("E4").NumberFormat = "d-mmm"

$ws6.Range("A5").Value = "PSG"
$ws6.Range("B5").Value = "Manchester City"
$ws6.Range("C5").Value = 2
$ws6.Range("D5").Value = 0
$ws6.Range("E5").Value = 44467
$ws6.Range("E5").NumberFormat = "d-mmm"

$ws6.Range("A6").Value = "Club Brugge"
$ws6.Range("B6").Value = "Manchester City"
$ws6.Range("C6").Value = 1
$ws6.Range("D6").Value = 5
$ws6.Range("E6").Value = 44488
$ws6.Range("E6").NumberFormat = "d-mmm"

$ws6.Range("A7").Value = "PSG"
$ws6.Range("B7").Value = "RB Leipzig"
$ws6.Range("C7").Value = 3
$ws6.Range("D7").Value = 2
$ws6.Range("E7").Value = 44488
$ws6.Range("E7").NumberFormat = "d-mmm"

$ws6.Range("A2").Select()

# ---------------------------------------------------------------------------
# Sheet 4: establishes the "bold + red + centered" divider style.
# ---------------------------------------------------------------------------
$ws4.Range("A1").Value = "team"
$ws4.Range("B1").Value = "plays"
$ws4.Range("C1").Value = "points"
$ws4.Range("A1:C1").Font.Bold = $true

$ws4.Range("A2").Value = "Group A"
$ws4.Range("A2:C2").Merge()
$ws4.Range("A2:C2").Font.Bold = $true
$ws4.Range("A2:C2").Font.Color = 255
$ws4.Range("A2:C2").HorizontalAlignment = -4108

$ws4.Range("A3").Value = "PSG"
$ws4.Range("B3").Value = 3
$ws4.Range("C3").Value = 7

$ws4.Range("A4").Value = "Manchester City"
$ws4.Range("B4").Value = 3
$ws4.Range("C4").Value = 6

$ws4.Range("A5").Value = "Club Brugge"
$ws4.Range("B5").Value = 3
$ws4.Range("C5").Value = 4

$ws4.Range("A6").Value = "RB Leipzig"
$ws4.Range("B6").Value = 3
$ws4.Range("C6").Value = 0

$ws4.Range("A7").Value = "Group B"
$ws4.Range("A7:C7").Merge()
$ws4.Range("A7:C7").Font.Bold = $true
$ws4.Range("A7:C7").Font.Color = 255
$ws4.Range("A7:C7").HorizontalAlignment = -4108

$ws4.Range("A8").Value = "Liverpool"
$ws4.Range("B8").Value = 3
$ws4.Range("C8").Value = 9

$ws4.Range("A9").Value = "Atletico Madrid"
$ws4.Range("B9").Value = 3
$ws4.Range("C9").Value = 4

$ws4.Range("A10").Value = "FC Porto"
$ws4.Range("B10").Value = 3
$ws4.Range("C10").Value = 4

$ws4.Range("A11").Value = "AC Milan"
$ws4.Range("B11").Value = 3
$ws4.Range("C11").Value = 0

$ws4.Range("B15").Select()

# ---------------------------------------------------------------------------
# Sheet 5: establishes the "centered + vcentered" merged-score-column style.
# ---------------------------------------------------------------------------
$ws5.Columns.Item(1).ColumnWidth = 18.1

$ws5.Range("A1").Value = "team"
$ws5.Range("B1").Value = "plays"
$ws5.Range("C1").Value = "points"
$ws5.Range("A1:C1").Font.Bold = $true

$ws5.Range("A2").Value = "PSG"
$ws5.Range("B2").Value = 3
$ws5.Range("C2").Value = 7

$ws5.Range("A3").Value = "Manchester City"
$ws5.Range("C3").Value = 6

$ws5.Range("A4").Value = "Club Brugge"
$ws5.Range("C4").Value = 4

$ws5.Range("A5").Value = "RB Leipzig"
$ws5.Range("C5").Value = 0

$ws5.Range("B2:B5").Merge()
$ws5.Range("B2:B5").HorizontalAlignment = -4108
$ws5.Range("B2:B5").VerticalAlignment = -4108

$ws5.Range("B6").Select()

# ---------------------------------------------------------------------------
# Update view state on the pre-existing sheets
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("multiple tables")
$ws2.Select()
$ws2.Range("A9:E15").Select()

$ws3 = $wb.Worksheets.Item("multiple tables with anchors")
$ws3.Select()
$ws3.Range("A2:C6").Select()

# the newly added "merged cells header" sheet ends up the active/selected tab
$ws6.Select()
